# Slide 11 ("Catch Curve Analysis in R") gets hidden from the slide show
# and a slow 2-second slide transition is added (PowerPoint 2010 timing
# extension p14:dur plus the legacy spd="slow" fallback).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)

# <p:sld ... show="0">  -- hide the slide from the slideshow
$s.SlideShowTransition.Hidden = $true

# <p:transition spd="slow" p14:dur="2000"/> (wrapped by PowerPoint in
# mc:AlternateContent / mc:Fallback for backward compatibility)
$s.SlideShowTransition.Duration = 2   # seconds -> p14:dur="2000" (ms)
$s.SlideShowTransition.Speed = 1      # ppTransitionSpeedSlow -> spd="slow"
